$d = $word.ActiveDocument

$replacements = @(
    @{old="702÷7=100, 2"; new="159÷5=31, 4"},
    @{old="808÷4=202, 0"; new="464÷8=58, 0"},
    @{old="333÷9=37, 0";  new="687÷3=229, 0"},
    @{old="821÷9=91, 2";  new="573÷7=81, 6"},
    @{old="379÷8=47, 3";  new="334÷3=111, 1"},
    @{old="959÷2=479, 1"; new="646÷3=215, 1"},
    @{old="396÷2=198, 0"; new="878÷8=109, 6"},
    @{old="314÷2=157, 0"; new="823÷8=102, 7"},
    @{old="731÷3=243, 2"; new="666÷4=166, 2"},
    @{old="443÷7=63, 2";  new="366÷8=45, 6"},
    @{old="629÷7=89, 6";  new="828÷6=138, 0"},
    @{old="900÷6=150, 0"; new="209÷4=52, 1"},
    @{old="702÷8=87, 6";  new="971÷2=485, 1"},
    @{old="509÷5=101, 4"; new="161÷9=17, 8"},
    @{old="582÷6=97, 0";  new="982÷7=140, 2"},
    @{old="699÷2=349, 1"; new="589÷7=84, 1"},
    @{old="242÷6=40, 2";  new="158÷6=26, 2"},
    @{old="810÷6=135, 0"; new="705÷8=88, 1"},
    @{old="872÷9=96, 8";  new="490÷7=70, 0"},
    @{old="165÷5=33, 0";  new="707÷7=101, 0"},
    @{old="283÷2=141, 1"; new="141÷5=28, 1"},
    @{old="668÷5=133, 3"; new="923÷5=184, 3"},
    @{old="188÷2=94, 0";  new="686÷3=228, 2"},
    @{old="684÷2=342, 0"; new="676÷5=135, 1"},
    @{old="472÷9=52, 4";  new="627÷9=69, 6"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
